$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.624.91"
$ws.Range("E2").Value = "  +3.56%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.608.03"
$ws.Range("E3").Value = "  +2.86%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'212.72"
$ws.Range("E5").Value = "  +1.24%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.521"
$ws.Range("E6").Value = "  +2.71%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - Solana
$ws.Range("D8").Value = "'26.86"
$ws.Range("E8").Value = "  +7.99%  "

# Row 9 - OKB
$ws.Range("D9").Value = "'43.59"
$ws.Range("E9").Value = "  -1.32%  "

# Row 10 - Cardano
$ws.Range("D10").Value = "'0.252"
$ws.Range("E10").Value = "  +2.63%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +2.61%  "

# Row 12 - TRON
$ws.Range("D12").Value = "'0.0908"

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "1.837.54"
$ws.Range("E13").Value = "  +2.80%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.608.76"
$ws.Range("E14").Value = "  +2.85%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "29.624.70"
$ws.Range("E15").Value = "  +3.39%  "

# Row 16 - Polygon
$ws.Range("E16").Value = "  +3.98%  "

# Row 17 - Polkadot
$ws.Range("E17").Value = "  +2.46%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "'63.44"
$ws.Range("E18").Value = "  +3.30%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'240.63"
$ws.Range("E19").Value = "  +6.06%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +3.87%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0693"
$ws.Range("E21").Value = "  +1.98%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.16%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +1.97%  "

# Row 24 - Avalanche
$ws.Range("D24").Value = "'9.24"
$ws.Range("E24").Value = "  +2.37%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  +0.72%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'154.47"
$ws.Range("E26").Value = "  +1.97%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  +3.64%  "

# Row 28 - Stellar
$ws.Range("D28").Value = "'0.109"
$ws.Range("E28").Value = "  +2.82%  "

# Row 29 - Cosmos
$ws.Range("D29").Value = "'6.42"
$ws.Range("E29").Value = "  +3.22%  "

# Row 30 - BinanceUSD
$ws.Range("E30").Value = "  -0.01%  "

# Row 31 - Hedera
$ws.Range("D31").Value = "'0.0474"
$ws.Range("E31").Value = "  +3.57%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +0.92%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  +1.55%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").Value = "'3.12"
$ws.Range("E34").Value = "  +4.47%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.411.61"
$ws.Range("E35").Value = "  +0.90%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  -0.36%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  +5.06%  "

# Row 38 - MXToken
$ws.Range("E38").Value = "  +5.39%  "

# Row 39 - HuobiToken
$ws.Range("E39").Value = "  +0.36%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +2.32%  "

# Row 41 - ImmutableX
$ws.Range("E41").Value = "  +4.41%  "

# Row 42 - RenderToken
$ws.Range("E42").Value = "  +2.56%  "

# Row 43 - Kaspa
$ws.Range("D43").Value = "'0.0490"
$ws.Range("E43").Value = "  +6.52%  "

# Row 44 - BitcoinSV
$ws.Range("E44").Value = "  +27.96%  "

# Row 45 - ARBITRUM
$ws.Range("D45").Value = "'0.799"
$ws.Range("E45").Value = "  +4.45%  "

# Row 46 - PaxDollar
$ws.Range("E46").Value = "  +0.00%  "

# Row 47 - Aave
$ws.Range("D47").Value = "'66.05"
$ws.Range("E47").Value = "  +3.46%  "

# Row 48 - FraxShare
$ws.Range("E48").Value = "  +1.33%  "

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = "1.747.72"
$ws.Range("E49").Value = "  +3.01%  "

# Row 50 - WEMIXToken
$ws.Range("D50").Value = "'0.870"
$ws.Range("E50").Value = "  +0.24%  "

# Row 51 - Quant
$ws.Range("D51").Value = "'86.62"
$ws.Range("E51").Value = "  +2.14%  "
